$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 499).
$ws.Range("C2:C499").Value = 45190
